# Updated PCB generated content
# Applies the regenerated Altium "Bill of Materials" report values to the
# "Part List Report" worksheet (rows 9-33 of the BOM table).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Part List Report")

# --- Header row 9: "Supplier Part Number 1" column label stays the same text ---
$ws.Range("J9").Value = "Supplier Part Number 1"

# --- Row 10: switch unit price / subtotal update ---
$ws.Range("M10").Value = 0.225
$ws.Range("N10").Value = 450

# --- Row 11: Digi-Key 220QBK-ND, stock depleted ---
$ws.Range("L11").Value = 0

# --- Row 12: quantity doubled (1 -> 2 pcs used), order qty doubled, stock depleted, subtotal doubled ---
$ws.Range("H12").Value = 2
$ws.Range("J12").Value = "220QBK-ND"
$ws.Range("K12").Value = 2000
$ws.Range("L12").Value = 0
$ws.Range("N12").Value = 22.88

# --- Row 12: component description now also covers the 1kOhm variant ---
$ws.Range("G12").Value = "General Type Carbon Film Resistor 3.3kOhm 1/4W 5% Axial Bulk, General Type Carbon Film Resistor 1kOhm 1/4W 5% Axial Bulk"

# --- Row 13: Digi-Key 220QBK-ND, stock depleted ---
$ws.Range("L13").Value = 0

# --- Row 17: LED array supplier stock/price/subtotal update ---
$ws.Range("L17").Value = 130
$ws.Range("M17").Value = 1.01
$ws.Range("N17").Value = 1008.17

# --- Row 18: LED red diffused, supplier stock update ---
$ws.Range("L18").Value = 67296

# --- Row 19: crystal, supplier stock update ---
$ws.Range("L19").Value = 908

# --- Row 20: connector, supplier stock update ---
$ws.Range("L20").Value = 123778

# --- Row 23: ATMEGA328P-PU now sourced from "Rs" instead of Farnell, new part #, new price ---
$ws.Range("I23").Value = "Rs"
$ws.Range("J23").Value = 1310276
$ws.Range("M23").Value = 2.75
$ws.Range("N23").Value = 2753.98

# --- Row 24: SN74HC595N now sourced from "Rs" (was Newark), new numeric part #, stock/price/subtotal update ---
$ws.Range("I24").Value = "Rs"
$ws.Range("J24").Value = 1000763
$ws.Range("L24").Value = 3750
$ws.Range("M24").Value = 0.45617
$ws.Range("N24").Value = 456.17

# --- Row 25: SN74LS47N, supplier stock update ---
$ws.Range("L25").Value = 5948

# --- Row 26: capacitor, unit price / subtotal update ---
$ws.Range("M26").Value = 0.36
$ws.Range("N26").Value = 720

# --- Row 27: capacitor, quantity used 2 -> 3, order qty 2000 -> 3000, unit price / subtotal update ---
$ws.Range("H27").Value = 3
$ws.Range("K27").Value = 3000
$ws.Range("M27").Value = 0.14195
$ws.Range("N27").Value = 425.84

# --- Row 29: SN74LS165AN now from Newark (shares index shift only), stock depleted, price/subtotal update ---
$ws.Range("I29").Value = "Newark"
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0.686
$ws.Range("N29").Value = 2058

# Totals in rows 30, 32, 33 are formulas (SUM / derived) and recalc
# automatically from the updated data above:
#   H30 =SUM(H10:H29), K30 =SUM(K10:K29), N30 =SUM(N10:N29)
#   L32 =N30, L33 =L32/H32
# D8/E8 (TODAY()/NOW()) also recalc automatically from the live clock.

# --- "Project Information" sheet: regenerated report metadata ---
$ws2 = $wb.Worksheets.Item("Project Information")
$ws2.Range("B7").Value = "53"
$ws2.Range("B8").Value = "1321h"
$ws2.Range("B9").Value = "08 Feb 2023"
$ws2.Range("B10").Value = "08 Feb 2023 1321h"
